$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 328.80768
$ws.Range("I33").Value = 267.3913
$ws.Range("J33").Value = 799.6667
$ws.Range("K33").Value = 267.3913
$ws.Range("L33").Value = 799.6667
$ws.Range("M33").Value = -38.3913
$ws.Range("N33").Value = -1257.6667
$ws.Range("H41").Value = 477.8
$ws.Range("I41").Value = 220
$ws.Range("J41").Value = 735.6
$ws.Range("K41").Value = 220
$ws.Range("L41").Value = 735.6
$ws.Range("M41").Value = 220
$ws.Range("N41").Value = -1615.6
$ws.Range("H53").Value = 311.35715
$ws.Range("I53").Value = 296.7143
$ws.Range("J53").Value = 326
$ws.Range("K53").Value = 296.7143
$ws.Range("L53").Value = 326
$ws.Range("M53").Value = 340.2857
$ws.Range("N53").Value = -1600
$ws.Range("H112").Value = 3473263.2
$ws.Range("I112").Value = 41667500
$ws.Range("J112").Value = 1059.9697
$ws.Range("K112").Value = 125002500
$ws.Range("L112").Value = 3179.9091
$ws.Range("M112").Value = -125001392
$ws.Range("N112").Value = -5395.909100000001
$ws.Range("H129").Value = 878.6
$ws.Range("I129").Value = 591.25
$ws.Range("J129").Value = 897.4426
$ws.Range("K129").Value = 1773.75
$ws.Range("L129").Value = 2692.3278
$ws.Range("M129").Value = 3226.25
$ws.Range("N129").Value = -12692.3278
$ws.Range("H137").Value = 1336.1111
$ws.Range("I137").Value = 996.1539
$ws.Range("J137").Value = 2220
$ws.Range("K137").Value = 2988.4617
$ws.Range("L137").Value = 6660
$ws.Range("M137").Value = -438.4616999999998
$ws.Range("N137").Value = -11760
$ws.Range("H138").Value = 3383.2632
$ws.Range("J138").Value = 3684.2354
$ws.Range("L138").Value = 11052.7062
$ws.Range("N138").Value = -21332.7062
$ws.Range("H141").Value = 2413.158
$ws.Range("J141").Value = 5184.1665
$ws.Range("L141").Value = 15552.4995
$ws.Range("N141").Value = -25912.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2601.111
$ws.Range("I2").Value = 1474.25
$ws.Range("J2").Value = 3502.6
$ws.Range("K2").Value = 1474.25
$ws.Range("L2").Value = 3502.6
$ws.Range("M2").Value = -1361.25
$ws.Range("N2").Value = -3728.6
$ws.Range("H32").Value = 3029.9348
$ws.Range("I32").Value = 2120.027
$ws.Range("K32").Value = 2120.027
$ws.Range("M32").Value = -1833.027
$ws.Range("H45").Value = 1925.027
$ws.Range("I45").Value = 1459.1818
$ws.Range("J45").Value = 2608.2666
$ws.Range("K45").Value = 1459.1818
$ws.Range("L45").Value = 2608.2666
$ws.Range("M45").Value = -1082.1818
$ws.Range("N45").Value = -3362.2666
$ws.Range("H116").Value = 2601.111
$ws.Range("I116").Value = 1474.25
$ws.Range("J116").Value = 3502.6
$ws.Range("K116").Value = 1474.25
$ws.Range("L116").Value = 3502.6
$ws.Range("M116").Value = 819.75
$ws.Range("N116").Value = -8090.6
$ws.Range("H132").Value = 21283.346
$ws.Range("I132").Value = 1538.7333
$ws.Range("J132").Value = 48207.816
$ws.Range("K132").Value = 4616.199900000001
$ws.Range("L132").Value = 144623.448
$ws.Range("M132").Value = -2086.199900000001
$ws.Range("N132").Value = -149683.448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2601.111
$ws.Range("I3").Value = 1474.25
$ws.Range("J3").Value = 3502.6
$ws.Range("K3").Value = 1474.25
$ws.Range("L3").Value = 3502.6
$ws.Range("M3").Value = -1360.25
$ws.Range("N3").Value = -3730.6
$ws.Range("H94").Value = 3907.111
$ws.Range("I94").Value = 1649.6666
$ws.Range("J94").Value = 5035.8335
$ws.Range("K94").Value = 1649.6666
$ws.Range("L94").Value = 5035.8335
$ws.Range("M94").Value = -1198.6666
$ws.Range("N94").Value = -5937.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 559.7778
$ws.Range("I22").Value = 609.6
$ws.Range("J22").Value = 497.5
$ws.Range("K22").Value = 609.6
$ws.Range("L22").Value = 497.5
$ws.Range("M22").Value = -259.6
$ws.Range("N22").Value = -1197.5
$ws.Range("H31").Value = 10252.941
$ws.Range("I31").Value = 12007.407
$ws.Range("J31").Value = 3485.7144
$ws.Range("K31").Value = 12007.407
$ws.Range("L31").Value = 3485.7144
$ws.Range("M31").Value = -11712.407
$ws.Range("N31").Value = -4075.7144
$ws.Range("H34").Value = 10252.941
$ws.Range("I34").Value = 12007.407
$ws.Range("J34").Value = 3485.7144
$ws.Range("K34").Value = 12007.407
$ws.Range("L34").Value = 3485.7144
$ws.Range("M34").Value = -11805.407
$ws.Range("N34").Value = -3889.7144
$ws.Range("H105").Value = 15626520
$ws.Range("I105").Value = 25000992
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 25000992
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -24999245
$ws.Range("N105").Value = -5894
$ws.Range("H132").Value = 17987.719
$ws.Range("I132").Value = 22565.262
$ws.Range("J132").Value = 6289.5557
$ws.Range("K132").Value = 67695.78599999999
$ws.Range("L132").Value = 18868.6671
$ws.Range("M132").Value = -65165.78599999999
$ws.Range("N132").Value = -23928.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 3250
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3250
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 9750
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -10062
$ws.Range("H64").Value = 2251.1428
$ws.Range("I64").Value = 1500
$ws.Range("K64").Value = 4500
$ws.Range("M64").Value = -4230
$ws.Range("H67").Value = 2251.1428
$ws.Range("I67").Value = 1500
$ws.Range("K67").Value = 4500
$ws.Range("M67").Value = -3564
$ws.Range("H131").Value = 815.02
$ws.Range("I131").Value = 440
$ws.Range("J131").Value = 834.7578999999999
$ws.Range("K131").Value = 1320
$ws.Range("L131").Value = 2504.2737
$ws.Range("M131").Value = 3720
$ws.Range("N131").Value = -12584.2737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1596.3334
$ws.Range("I107").Value = 428
$ws.Range("K107").Value = 428
$ws.Range("M107").Value = 1492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6453732.5
$ws.Range("I136").Value = 16130032
$ws.Range("K136").Value = 48390096
$ws.Range("M136").Value = -48387546
